# "Se agrega cuadro para agregar al checklist"
#
# Adds a new "ID" checklist column in front of the existing data on the
# "KG" sheet: a new column A is inserted (pushing Producto/Precio/g/len
# from A:D to B:E), headed "ID", with the first product row tagged
# "001s" and every other row tagged "00".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KG")

# Insert a blank column at A; existing A:D (Producto/Precio/g/len) and
# their formulas shift right to B:E automatically.
$ws.Columns.Item(1).Insert()

# Header for the new checklist column.
$ws.Range("A1").Value2 = "ID"

# Make sure the ID values are stored as text (they look numeric, e.g.
# "00", so force a text number format before writing them).
$ws.Range("A2:A91").NumberFormat = "@"

# Row 2 (first product) is marked "001s"; every remaining product row
# is marked "00". Fill the plain "00" rows first so the shared-string
# table order matches (00 before 001s).
for ($r = 3; $r -le 91; $r++) {
    $ws.Cells.Item($r, 1).Value2 = "00"
}
$ws.Range("A2").Value2 = "001s"

# The AutoFilter should keep covering only the original data columns
# (now B:E), not the new ID column.
$ws.AutoFilterMode = $false
$null = $ws.Range("B1:E91").AutoFilter()

# Keep the hidden _FilterDatabase defined name for KG in sync with the
# shifted data range.
$fdb = $wb.Names.Item("KG!_FilterDatabase")
$fdb.RefersTo = "=KG!`$B`$1:`$E`$91"

# Match the author's final selection/cursor position.
$ws.Activate()
$null = $ws.Range("A3").Select()
